# Updated cryptos list on Sat Jun 29 11:08:47 UTC 2024 with GitHub Actions
# Refresh price / 1h-volume figures (and swap the EnergySwap/Hedera row
# order) to match the latest scrape.
#
# Note: several "Price" (column D) values look like plain numbers
# (e.g. "1.00", "28.05") but must stay stored as literal text, matching
# the rest of the column (t="inlineStr" in the original workbook). A
# leading apostrophe forces Excel to keep the exact textual
# representation (trailing zeros etc.) instead of coercing it to a
# Double and losing formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.936.22'
$ws.Range("E2").Value = '  -0.91%  '

$ws.Range("D3").Value = '3.394.62'
$ws.Range("E3").Value = '  -1.45%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '''572.70'

$ws.Range("D6").Value = '''142.51'
$ws.Range("E6").Value = '  -1.88%  '

$ws.Range("D7").Value = '3.395.24'
$ws.Range("E7").Value = '  -1.46%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("E9").Value = '  -0.79%  '

$ws.Range("E12").Value = '  +1.42%  '

$ws.Range("D13").Value = '3.974.41'
$ws.Range("E13").Value = '  -1.47%  '

$ws.Range("E14").Value = '  +2.09%  '

$ws.Range("D15").Value = '''28.05'
$ws.Range("E15").Value = '  -1.52%  '

$ws.Range("E16").Value = '  -1.28%  '

$ws.Range("D17").Value = '3.399.54'
$ws.Range("E17").Value = '  -1.45%  '

$ws.Range("D18").Value = '61.007.98'
$ws.Range("E18").Value = '  -1.04%  '

$ws.Range("D19").Value = '''6.13'
$ws.Range("E19").Value = '  -3.51%  '

$ws.Range("E20").Value = '  -3.16%  '

$ws.Range("D21").Value = '''8.95'
$ws.Range("E21").Value = '  -4.92%  '

$ws.Range("D22").Value = '''383.15'
$ws.Range("E22").Value = '  -4.76%  '

$ws.Range("D23").Value = '''0.558'
$ws.Range("E23").Value = '  -1.81%  '

$ws.Range("D24").Value = '''74.32'
$ws.Range("E24").Value = '  -0.14%  '

$ws.Range("D25").Value = '''0.998'
$ws.Range("E25").Value = '  -0.29%  '

$ws.Range("D26").Value = '''0.0000117'
$ws.Range("E26").Value = '  -5.03%  '

$ws.Range("D27").Value = '3.533.56'
$ws.Range("E27").Value = '  -1.57%  '

$ws.Range("E28").Value = '  -1.31%  '

$ws.Range("E29").Value = '  -0.12%  '

$ws.Range("D30").Value = '''7.38'
$ws.Range("E30").Value = '  -3.19%  '

$ws.Range("D31").Value = '''8.01'
$ws.Range("E31").Value = '  -3.01%  '

$ws.Range("E32").Value = '  -1.19%  '

$ws.Range("D33").Value = '''1.41'
$ws.Range("E33").Value = '  -2.42%  '

$ws.Range("E34").Value = '  -0.06%  '

$ws.Range("D35").Value = '''23.50'
$ws.Range("E35").Value = '  -1.72%  '

$ws.Range("E36").Value = '  -0.75%  '

$ws.Range("D37").Value = '''167.62'
$ws.Range("E37").Value = '  +0.32%  '

$ws.Range("D38").Value = '3.425.75'
$ws.Range("E38").Value = '  -1.36%  '

$ws.Range("E39").Value = '  -2.74%  '

$ws.Range("E40").Value = '  -4.90%  '

$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = '''27.62'
$ws.Range("E41").Value = '  +1.22%  '

$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D42").Value = '''0.0772'
$ws.Range("E42").Value = '  -2.41%  '

$ws.Range("D43").Value = '''0.781'
$ws.Range("E43").Value = '  -2.51%  '

$ws.Range("E44").Value = '  +0.00%  '

$ws.Range("D45").Value = '''4.43'
$ws.Range("E45").Value = '  -2.14%  '

$ws.Range("E46").Value = '  -3.65%  '

$ws.Range("E47").Value = '  -1.25%  '

$ws.Range("D48").Value = '2.480.88'
$ws.Range("E48").Value = '  -5.07%  '

$ws.Range("D49").Value = '''6.81'

$ws.Range("D50").Value = '''22.99'
$ws.Range("E50").Value = '  -0.57%  '

$ws.Range("E51").Value = '  +0.79%  '
